$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E ("nr tel"), shifting "nr tel" and "notatka" right
$ws.Columns("E:E").Insert()

# New header for inserted column
$ws.Range("E1").Value = "grupa"

# New data values for the inserted column
$ws.Range("E2").Value = "jun1"
$ws.Range("E3").Value = "ważna grupa2 AĄŁ"

# Rename "notatka" header (now in column G) to "notatka rekrutacyjna"
$ws.Range("G1").Value = "notatka rekrutacyjna"

# Update selection to match target (G5)
$ws.Range("G5").Select()
